# Auto: Weekly update of data
# Appends 4 new weekly match rows (103-106) to the "Main" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$data = @(
    @{ row = 103; A = 862; B = 311; C = 311; D = 18997; E = 20534;          G = 17356; H = 17238; I = 17882; J = 17772; K = 18153; L = 15984; M = 17725; N = -8; O = 10; P = 3;  Q = -7;  R = -5; S = 1; T = -2 },
    @{ row = 104; A = 863; B = 311; C = 221; D = 18026; E = 20462; F = 13275; G = 17705; H = 17511; I = 18248; J = 17669; K = 17520; L = 17230; M = 17335; N = 5;  O = 3;  P = 0;  Q = -13; R = 0;  S = 6; T = -1 },
    @{ row = 105; A = 864; B = 32;  C = 311; D = 16474; E = 15963; F = 13632; G = 17477; H = 18064; I = 17508; J = 16898; K = 12021; L = 14101; M = 19541; N = 0;  O = 0;  P = -1; Q = -8;  R = 7;  S = 4; T = -2 },
    @{ row = 106; A = 865; B = 311; C = 221; D = 13267; E = 17098; F = 17725; G = 14999; H = 14813; I = 19395; J = 10415; K = 16483; L = 14477;            N = 7;  O = 7;  P = -4; Q = 3;   R = -1; S = -1; T = -2 }
)

foreach ($rowData in $data) {
    $r = $rowData.row
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")) {
        if ($rowData.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $rowData[$col]
        }
    }
}

$ws.Range("N106").Select()
